# NATMI ligand-receptor pair workbook was regenerated with updated TPM values.
# The "Resolving-Mac" target cluster is no longer present among the target
# clusters (the 4th target per sending cluster), so every row whose "Target
# cluster" was "Resolving-Mac" is dropped, and all derived/specificity
# columns (ligand/receptor detection + specificity + edge weights) are
# recomputed against the remaining 3 target clusters. Net effect: the sheet
# shrinks from 12 data rows (A1:T13) to 9 data rows (A1:T10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer present (old rows 11-13) so the sheet shrinks to A1:T10
$ws.Range("A11:T13").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2412233333333333
$ws.Range("H2").Value = 0.72367
$ws.Range("I2").Value = 0.2038054651530871
$ws.Range("J2").Value = 0.2038054651530872
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.227228
$ws.Range("N2").Value = 0.681684
$ws.Range("O2").Value = 0.2376267857721762
$ws.Range("P2").Value = 0.2376267857721762
$ws.Range("Q2").Value = 0.05481269558666667
$ws.Range("R2").Value = 0.49331426028
$ws.Range("S2").Value = 0.04842963760713136
$ws.Range("T2").Value = 0.04842963760713136

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2412233333333333
$ws.Range("H3").Value = 0.72367
$ws.Range("I3").Value = 0.2038054651530871
$ws.Range("J3").Value = 0.2038054651530872
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6537306666666667
$ws.Range("N3").Value = 1.961192
$ws.Range("O3").Value = 0.6836477770376096
$ws.Range("P3").Value = 0.6836477770376095
$ws.Range("Q3").Value = 0.1576950905155556
$ws.Range("R3").Value = 1.41925581464
$ws.Range("S3").Value = 0.139331153200024
$ws.Range("T3").Value = 0.139331153200024

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2412233333333333
$ws.Range("H4").Value = 0.72367
$ws.Range("I4").Value = 0.2038054651530871
$ws.Range("J4").Value = 0.2038054651530872
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07528033333333332
$ws.Range("N4").Value = 0.225841
$ws.Range("O4").Value = 0.0787254371902143
$ws.Range("P4").Value = 0.0787254371902143
$ws.Range("Q4").Value = 0.01815937294111111
$ws.Range("R4").Value = 0.16343435647
$ws.Range("S4").Value = 0.01604467434593177
$ws.Range("T4").Value = 0.01604467434593177

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.133983
$ws.Range("H5").Value = 0.401949
$ws.Range("I5").Value = 0.1131999432238703
$ws.Range("J5").Value = 0.1131999432238703
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.227228
$ws.Range("N5").Value = 0.681684
$ws.Range("O5").Value = 0.2376267857721762
$ws.Range("P5").Value = 0.2376267857721762
$ws.Range("Q5").Value = 0.03044468912399999
$ws.Range("R5").Value = 0.274002202116
$ws.Range("S5").Value = 0.02689933865788113
$ws.Range("T5").Value = 0.02689933865788113

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.133983
$ws.Range("H6").Value = 0.401949
$ws.Range("I6").Value = 0.1131999432238703
$ws.Range("J6").Value = 0.1131999432238703
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6537306666666667
$ws.Range("N6").Value = 1.961192
$ws.Range("O6").Value = 0.6836477770376096
$ws.Range("P6").Value = 0.6836477770376095
$ws.Range("Q6").Value = 0.087588795912
$ws.Range("R6").Value = 0.788299163208
$ws.Range("S6").Value = 0.07738888954578255
$ws.Range("T6").Value = 0.07738888954578255

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.133983
$ws.Range("H7").Value = 0.401949
$ws.Range("I7").Value = 0.1131999432238703
$ws.Range("J7").Value = 0.1131999432238703
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07528033333333332
$ws.Range("N7").Value = 0.225841
$ws.Range("O7").Value = 0.0787254371902143
$ws.Range("P7").Value = 0.0787254371902143
$ws.Range("Q7").Value = 0.010086284901
$ws.Range("R7").Value = 0.09077656410899999
$ws.Range("S7").Value = 0.008911715020206625
$ws.Range("T7").Value = 0.008911715020206627

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cxcl13"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8083896666666667
$ws.Range("H8").Value = 2.425169
$ws.Range("I8").Value = 0.6829945916230425
$ws.Range("J8").Value = 0.6829945916230425
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.227228
$ws.Range("N8").Value = 0.681684
$ws.Range("O8").Value = 0.2376267857721762
$ws.Range("P8").Value = 0.2376267857721762
$ws.Range("Q8").Value = 0.1836887671773333
$ws.Range("R8").Value = 1.653198904596
$ws.Range("S8").Value = 0.1622978095071637
$ws.Range("T8").Value = 0.1622978095071637

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cxcl13"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8083896666666667
$ws.Range("H9").Value = 2.425169
$ws.Range("I9").Value = 0.6829945916230425
$ws.Range("J9").Value = 0.6829945916230425
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6537306666666667
$ws.Range("N9").Value = 1.961192
$ws.Range("O9").Value = 0.6836477770376096
$ws.Range("P9").Value = 0.6836477770376095
$ws.Range("Q9").Value = 0.5284691157164445
$ws.Range("R9").Value = 4.756222041448
$ws.Range("S9").Value = 0.466927734291803
$ws.Range("T9").Value = 0.466927734291803

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Cxcl13"
$ws.Range("C10").Value = "Ackr4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8083896666666667
$ws.Range("H10").Value = 2.425169
$ws.Range("I10").Value = 0.6829945916230425
$ws.Range("J10").Value = 0.6829945916230425
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.07528033333333332
$ws.Range("N10").Value = 0.225841
$ws.Range("O10").Value = 0.0787254371902143
$ws.Range("P10").Value = 0.0787254371902143
$ws.Range("Q10").Value = 0.06085584356988888
$ws.Range("R10").Value = 0.547702592129
$ws.Range("S10").Value = 0.0537690478240759
$ws.Range("T10").Value = 0.0537690478240759
